{"js": "// Replace the 25 two-digit multiplication problems in the practice-sheet\n// table with the new values from the commit. Every \"old\" value below is\n// unique within the document, so an exact-text search-and-replace safely\n// targets the correct cell each time.\nconst pairs = [\n  [\"35\u00d794=\", \"38\u00d748=\"],\n  [\"98\u00d788=\", \"68\u00d785=\"],\n  [\"61\u00d791=\", \"46\u00d755=\"],\n  [\"92\u00d726=\", \"16\u00d784=\"],\n  [\"59\u00d749=\", \"64\u00d754=\"],\n  [\"19\u00d726=\", \"66\u00d755=\"],\n  [\"60\u00d740=\", \"37\u00d761=\"],\n  [\"35\u00d754=\", \"86\u00d781=\"],\n  [\"80\u00d790=\", \"49\u00d737=\"],\n  [\"39\u00d785=\", \"79\u00d753=\"],\n  [\"62\u00d739=\", \"13\u00d748=\"],\n  [\"31\u00d716=\", \"19\u00d761=\"],\n  [\"98\u00d757=\", \"22\u00d783=\"],\n  [\"62\u00d770=\", \"49\u00d766=\"],\n  [\"67\u00d740=\", \"26\u00d731=\"],\n  [\"93\u00d794=\", \"60\u00d747=\"],\n  [\"17\u00d722=\", \"94\u00d759=\"],\n  [\"23\u00d734=\", \"18\u00d750=\"],\n  [\"34\u00d719=\", \"36\u00d796=\"],\n  [\"54\u00d730=\", \"35\u00d795=\"],\n  [\"76\u00d786=\", \"14\u00d784=\"],\n  [\"43\u00d715=\", \"51\u00d723=\"],\n  [\"86\u00d779=\", \"20\u00d735=\"],\n  [\"36\u00d782=\", \"23\u00d731=\"],\n  [\"74\u00d788=\", \"52\u00d791=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit multiplication problems in the practice-sheet\n# table with the new values from the commit. Every \"old\" value below is\n# unique within the document, so an exact Find/Replace targets the\n# correct cell each time.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"35\u00d794=\", \"38\u00d748=\"),\n    @(\"98\u00d788=\", \"68\u00d785=\"),\n    @(\"61\u00d791=\", \"46\u00d755=\"),\n    @(\"92\u00d726=\", \"16\u00d784=\"),\n    @(\"59\u00d749=\", \"64\u00d754=\"),\n    @(\"19\u00d726=\", \"66\u00d755=\"),\n    @(\"60\u00d740=\", \"37\u00d761=\"),\n    @(\"35\u00d754=\", \"86\u00d781=\"),\n    @(\"80\u00d790=\", \"49\u00d737=\"),\n    @(\"39\u00d785=\", \"79\u00d753=\"),\n    @(\"62\u00d739=\", \"13\u00d748=\"),\n    @(\"31\u00d716=\", \"19\u00d761=\"),\n    @(\"98\u00d757=\", \"22\u00d783=\"),\n    @(\"62\u00d770=\", \"49\u00d766=\"),\n    @(\"67\u00d740=\", \"26\u00d731=\"),\n    @(\"93\u00d794=\", \"60\u00d747=\"),\n    @(\"17\u00d722=\", \"94\u00d759=\"),\n    @(\"23\u00d734=\", \"18\u00d750=\"),\n    @(\"34\u00d719=\", \"36\u00d796=\"),\n    @(\"54\u00d730=\", \"35\u00d795=\"),\n    @(\"76\u00d786=\", \"14\u00d784=\"),\n    @(\"43\u00d715=\", \"51\u00d723=\"),\n    @(\"86\u00d779=\", \"20\u00d735=\"),\n    @(\"36\u00d782=\", \"23\u00d731=\"),\n    @(\"74\u00d788=\", \"52\u00d791=\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]2)\n}\n"}
